$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, pushing the existing row 3 (and below) down to row 4.
$ws.Rows.Item(3).Insert()

# Copy the date-format style (used in column D for rows 2-3) to the new D3 cell.
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# Populate the new row 3 with the updated record.
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 44875
$ws.Cells.Item(3, 5).Value = 5
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100104
$ws.Cells.Item(3, 8).Value = "Frutos de pepita"
$ws.Cells.Item(3, 9).Value = 100104004
$ws.Cells.Item(3, 10).Value = "Níspero"
$ws.Cells.Item(3, 11).Value = "Golden Nugget"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 50
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 16000
$ws.Cells.Item(3, 16).Value = 16000
$ws.Cells.Item(3, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(3, 19).Value = 1600
$ws.Cells.Item(3, 20).Value = 10
